# Update question difficulty level
# - Moves the existing "recommended_level" column (F) to K
# - Inserts four new computed columns: frequency (F), frequency_occurrence (G),
#   frequency_occurrence_probab (H), max_probab (I), max_probab_percentage (J)
# - Keeps "recommended_level" as the header for the (new) last column K

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First, move the old F column (recommended_level header + values) over to K,
# since the diff shows the previous F data ending up unchanged in K.
$ws.Range("F1:F17").Copy()
$ws.Range("K1").PasteSpecial()

# Give the new header cells (G1:K1) the same bold/bordered header formatting
# used by the rest of row 1 (B1:F1), then fill in their captions.
$ws.Range("B1").Copy()
$ws.Range("G1:K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("F1").Value = "frequency"
$ws.Range("G1").Value = "frequency_occurrence"
$ws.Range("H1").Value = "frequency_occurrence_probab"
$ws.Range("I1").Value = "max_probab"
$ws.Range("J1").Value = "max_probab_percentage"
$ws.Range("K1").Value = "recommended_level"

# Row 2 has no computed values (matches source row that had blank recommended_level)
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""

# Per-row computed data: frequency, frequency_occurrence, frequency_occurrence_probab,
# max_probab, max_probab_percentage (recommended_level in K is already correct from the copy above)
$data = @{
    3  = @(8,  '{"L3":8}',          '{"L3":1.0}',               1,     "100.00")
    4  = @(8,  '{"L2":5,"L1":3}',   '{"L2":0.625,"L1":0.375}',  0.625, "62.50")
    5  = @(24, '{"L3":12,"L2":12}', '{"L3":0.5,"L2":0.5}',      0.5,   "50.00")
    6  = @(8,  '{"L1":8}',          '{"L1":1.0}',               1,     "100.00")
    7  = @(8,  '{"L3":8}',          '{"L3":1.0}',               1,     "100.00")
    8  = @(24, '{"L3":12,"L2":12}', '{"L3":0.5,"L2":0.5}',      0.5,   "50.00")
    9  = @(8,  '{"L3":8}',          '{"L3":1.0}',               1,     "100.00")
    10 = @(8,  '{"L2":6,"L1":2}',   '{"L2":0.75,"L1":0.25}',    0.75,  "75.00")
    11 = @(8,  '{"L1":4,"L2":4}',   '{"L1":0.5,"L2":0.5}',      0.5,   "50.00")
    12 = @(8,  '{"L3":5,"L2":3}',   '{"L3":0.625,"L2":0.375}',  0.625, "62.50")
    13 = @(8,  '{"L3":8}',          '{"L3":1.0}',               1,     "100.00")
    14 = @(8,  '{"L3":6,"L2":2}',   '{"L3":0.75,"L2":0.25}',    0.75,  "75.00")
    15 = @(8,  '{"L3":5,"L2":3}',   '{"L3":0.625,"L2":0.375}',  0.625, "62.50")
    16 = @(7,  '{"L1":7}',          '{"L1":1.0}',               1,     "100.00")
    17 = @(8,  '{"L3":8}',          '{"L3":1.0}',               1,     "100.00")
}

# The percentage column stores text like "100.00"/"62.50", so force a text
# number format before assigning or Excel would otherwise coerce it to a number.
$ws.Range("J3:J17").NumberFormat = "@"

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 6).Value  = $vals[0]   # F: frequency
    $ws.Cells.Item($row, 7).Value  = $vals[1]   # G: frequency_occurrence
    $ws.Cells.Item($row, 8).Value  = $vals[2]   # H: frequency_occurrence_probab
    $ws.Cells.Item($row, 9).Value  = $vals[3]   # I: max_probab
    $ws.Cells.Item($row, 10).Value = $vals[4]   # J: max_probab_percentage
}

$ws.Range("A1").Select()
